# Populate the "Definition" column (D) on the Concepts sheet with the
# same value as the "Display" column (C) for each concept row (2-7).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

for ($row = 2; $row -le 7; $row++) {
    $display = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 4).Value = $display
}
